# Restore D8 ("Integer max") on the Rules sheet from 11 to 19.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$ws.Range("D8").Value = 19
